$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.272.41'
$ws.Range("E2").Value = '  +2.71%  '

$ws.Range("D3").Value = '2.424.21'
$ws.Range("E3").Value = '  +2.12%  '

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '307.77'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.63%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '100.91'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +4.09%  '

$ws.Range("E7").Value = '  +1.67%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  +0.70%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '35.29'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +3.29%  '

$ws.Range("E11").Value = '  +2.06%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '18.89'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +3.17%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.124'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +2.16%  '

$ws.Range("E14").Value = '  +2.00%  '

$ws.Range("D15").Value = '2.803.70'
$ws.Range("E15").Value = '  +2.17%  '

$ws.Range("D16").Value = '2.433.45'
$ws.Range("E16").Value = '  +1.63%  '

$ws.Range("E17").Value = '  +3.45%  '

$ws.Range("D18").Value = '44.200.64'
$ws.Range("E18").Value = '  +2.61%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '12.31'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.43%  '

$ws.Range("E20").Value = '  +1.77%  '

$ws.Range("D21").Value = '0.0₃0907'
$ws.Range("E21").Value = '  +2.22%  '

$ws.Range("E22").Value = '  +0.36%  '

$ws.Range("E23").Value = '  +5.63%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '240.41'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +2.20%  '

$ws.Range("E25").Value = '  +1.65%  '

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.16%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '25.18'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.84%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.32'
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '9.59'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +5.16%  '

$ws.Range("E30").Value = '  +5.39%  '

$ws.Range("E31").Value = '  +12.43%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '18.77'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +7.88%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.18'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +2.24%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("E35").Value = '  +1.48%  '

$ws.Range("E36").Value = '  +3.60%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '131.02'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +25.10%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '4.46'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +4.47%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.89'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +3.68%  '

$ws.Range("E40").Value = '  -0.77%  '

$ws.Range("E41").Value = '  +0.77%  '

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '21.39'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -4.29%  '

$ws.Range("E43").Value = '  +2.68%  '

$ws.Range("D44").Value = '1.949.14'
$ws.Range("E44").Value = '  -0.29%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.87'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +4.85%  '

$ws.Range("E47").Value = '  +2.46%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.66'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +10.72%  '

$ws.Range("E49").Value = '  +1.69%  '

$ws.Range("E50").Value = '  +2.69%  '

$ws.Range("E51").Value = '  +0.97%  '

